{"js": "// The document's numbered \"Brief Acceptance of Harm Reduction Scales\" list\n// items were re-shuffled/re-worded, two brand-new items were added, and the\n// \"Scoring\" section at the bottom was updated to reference the new item\n// numbers (1-9 items instead of 1-8, reverse-coded items 11/13/14/16\n// instead of 9/11/12/14, and Principles subscale 10-16 instead of 9-15).\n//\n// We locate each affected paragraph by its position in\n// context.document.body.paragraphs (stable regardless of the run-split\n// inside it) and replace its whole text in one shot with\n// Word.InsertLocation.replace, which keeps the paragraph's existing run\n// formatting (font/size/color) intact.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph index -> new full paragraph text, for the 8 numbered\n// \"acceptance\" statements that changed (paragraph 12, item 9, is untouched\n// by the diff and is intentionally omitted here).\nconst itemReplacements = {\n  4: \"People who use drugs should have access to safe injection supplies (sterile needles and syringes).\",\n  5: \"People who inject drugs should be able to do so in a way that prevents them from causing further harm to their health.\",\n  6: \"People who use drugs should have access to tools to test what's in their drugs.\",\n  7: \"People who use drugs should have access to supervised places where they can consume drugs safely.\",\n  8: \"People who use drugs should have access to a legal, non-contaminated drug supply.\",\n  9: \"People should be able to use drugs safely.\",\n  10: \"Racism affects the health of people who use drugs.\",\n  11: \"People who seek medical assistance for overdoses should be protected from drug charges, arrests, and prosecutions.\",\n};\n\nfor (const key of Object.keys(itemReplacements)) {\n  const idx = parseInt(key, 10);\n  paragraphs.items[idx].insertText(itemReplacements[key], Word.InsertLocation.replace);\n}\n\n// Scoring section: update the reverse-coded item list and the two subscale\n// ranges so they match the renumbered items above.\nparagraphs.items[22].insertText(\n  \"Items 11, 13, 14, & 16 are reverse coded\",\n  Word.InsertLocation.replace\n);\nparagraphs.items[23].insertText(\n  \"Harm Reduction Strategies Subscale: 1-9\",\n  Word.InsertLocation.replace\n);\nparagraphs.items[24].insertText(\n  \"Harm Reduction Principles Subscale 10-16\",\n  Word.InsertLocation.replace\n);\n\nawait context.sync();\n", "ps1": "# The document's numbered \"Brief Acceptance of Harm Reduction Scales\" list\n# items were re-shuffled/re-worded, two brand-new items were added, and the\n# \"Scoring\" section at the bottom was updated to reference the new item\n# numbers (1-9 items instead of 1-8, reverse-coded items 11/13/14/16\n# instead of 9/11/12/14, and Principles subscale 10-16 instead of 9-15).\n#\n# Several of the affected paragraphs originally share identical text (e.g.\n# paragraphs 5 and 7 both read \"People who inject drugs...\"), so we cannot\n# safely use a document-wide Find/Replace. Instead we resolve each paragraph\n# by its fixed position (these edits only change text, never add/remove\n# paragraphs, so positions stay stable across the loop) and overwrite just\n# that paragraph's text (excluding its trailing paragraph mark, so the\n# paragraph count/structure is untouched). This also keeps the existing run\n# formatting (font/size/color) intact and avoids Word's automatic\n# straight-quote -> curly-quote substitution that Find/Replace would apply.\n\n$d = $word.ActiveDocument\n\nfunction Set-ParagraphText($doc, [int]$index, [string]$newText) {\n    $p = $doc.Paragraphs.Item($index)\n    $start = $p.Range.Start\n    $end = $p.Range.End - 1   # exclude the trailing paragraph mark\n    $r = $doc.Range($start, $end)\n    $r.Text = $newText\n}\n\n# The 8 numbered \"acceptance\" statements that changed (paragraph 13, item 9,\n# \"Possession of drug paraphernalia...\", is untouched by the diff and is\n# intentionally left out here).\nSet-ParagraphText $d 5  \"People who use drugs should have access to safe injection supplies (sterile needles and syringes).\"\nSet-ParagraphText $d 6  \"People who inject drugs should be able to do so in a way that prevents them from causing further harm to their health.\"\nSet-ParagraphText $d 7  \"People who use drugs should have access to tools to test what's in their drugs.\"\nSet-ParagraphText $d 8  \"People who use drugs should have access to supervised places where they can consume drugs safely.\"\nSet-ParagraphText $d 9  \"People who use drugs should have access to a legal, non-contaminated drug supply.\"\nSet-ParagraphText $d 10 \"People should be able to use drugs safely.\"\nSet-ParagraphText $d 11 \"Racism affects the health of people who use drugs.\"\nSet-ParagraphText $d 12 \"People who seek medical assistance for overdoses should be protected from drug charges, arrests, and prosecutions.\"\n\n# Scoring section: update the reverse-coded item list and the two subscale\n# ranges so they match the renumbered items above.\nSet-ParagraphText $d 23 \"Items 11, 13, 14, & 16 are reverse coded\"\nSet-ParagraphText $d 24 \"Harm Reduction Strategies Subscale: 1-9\"\nSet-ParagraphText $d 25 \"Harm Reduction Principles Subscale 10-16\"\n"}
